# "Added last minute updates"
#
# 1. First paragraph (the hidden **ID__...__ID** merge-field placeholder):
#      - rename the placeholder id from
#          **ID__AFFARS_pgi_5315_topic_15__ID**
#        to
#          **ID__AFFARS_AFMC_PGI_5315__ID**
#      - drop the trailing run that held only a single space character
#      - add a (spaceless/borderless-color) paragraph border on all four
#        sides with 5pt padding
#      - widen the left indent from 120 -> 225 twips (6pt -> 11.25pt)

$d = $word.ActiveDocument

$oldId = "**ID__AFFARS_pgi_5315_topic_15__ID**"
$newId = "**ID__AFFARS_AFMC_PGI_5315__ID**"

$p1 = $d.Paragraphs.Item(1)

# Rename the placeholder id. Restrict the Find to paragraph 1's own range so
# it can never spill over into the following heading paragraph.
$idRange = $p1.Range.Duplicate
[void]$idRange.Find.Execute($oldId, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $newId, 2)

# The paragraph used to carry a second run with nothing but a single space
# in it; that run is removed entirely by the edit, so drop the now-stray
# space left behind by the rename above.
$p1 = $d.Paragraphs.Item(1)
$spaceRange = $p1.Range.Duplicate
[void]$spaceRange.Find.Execute(" ", $true, $false, $false, $false, $false, `
                                $true, 1, $false, "", 2)

# Re-fetch the paragraph (text changed) and apply the new paragraph-level
# formatting: a thin border around the paragraph and a wider left indent.
$p1 = $d.Paragraphs.Item(1)
$p1.Format.LeftIndent = 11.25  # 225 twips

$borders = $p1.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5
